$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $value) {
    # Force the range to be treated as plain text so Excel does not
    # auto-coerce numeric-looking strings (e.g. "1.010", "10.00") into
    # numbers and strip the formatting-significant trailing zeros.
    # ClearFormats afterwards drops the temporary "@" number format so the
    # cell's style index reverts to the sheet default (unstyled), matching
    # the original workbook's cell styling.
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "26.529.30"
Set-TextValue $ws.Range("E2") "  -2.68%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.812.70"
Set-TextValue $ws.Range("E3") "  -2.43%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.007"
Set-TextValue $ws.Range("E4") "  +0.46%  "

# Row 5 - was USDC, now BNB
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue $ws.Range("D5") "308.77"
Set-TextValue $ws.Range("E5") "  -1.71%  "

# Row 6 - was BNB, now USDC
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D6") "1.007"
Set-TextValue $ws.Range("E6") "  +0.46%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.4572"
Set-TextValue $ws.Range("E7") "  -1.83%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.3668"
Set-TextValue $ws.Range("E8") "  -1.18%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.07130"
Set-TextValue $ws.Range("E9") "  -2.39%  "

# Row 10 - Polygon
Set-TextValue $ws.Range("D10") "0.8802"
Set-TextValue $ws.Range("E10") "  -1.38%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.07758"
Set-TextValue $ws.Range("E11") "  -1.54%  "

# Row 12 - Solana
Set-TextValue $ws.Range("D12") "19.36"
Set-TextValue $ws.Range("E12") "  -3.72%  "

# Row 13 - WrappedEther
Set-TextValue $ws.Range("D13") "1.832.64"
Set-TextValue $ws.Range("E13") "  -1.37%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "5.293"
Set-TextValue $ws.Range("E14") "  -2.12%  "

# Row 15 - Chainlink
Set-TextValue $ws.Range("D15") "6.376"
Set-TextValue $ws.Range("E15") "  -2.23%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "86.65"
Set-TextValue $ws.Range("E16") "  -5.34%  "

# Row 17 - BinanceUSD
Set-TextValue $ws.Range("D17") "1.010"
Set-TextValue $ws.Range("E17") "  +0.74%  "

# Row 18 - ShibaInu
Set-TextValue $ws.Range("D18") "0.000008595"
Set-TextValue $ws.Range("E18") "  -3.88%  "

# Row 19 - Dai
Set-TextValue $ws.Range("D19") "1.007"
Set-TextValue $ws.Range("E19") "  +0.54%  "

# Row 20 - WrappedBTC
Set-TextValue $ws.Range("D20") "26.587.36"
Set-TextValue $ws.Range("E20") "  -2.57%  "

# Row 21 - Avalanche
Set-TextValue $ws.Range("E21") "  -3.44%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("E22") "  -1.59%  "

# Row 23 - Cosmos
Set-TextValue $ws.Range("E23") "  -0.56%  "

# Row 24 - Toncoin
Set-TextValue $ws.Range("D24") "1.983"
Set-TextValue $ws.Range("E24") "  -3.78%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "151.40"
Set-TextValue $ws.Range("E25") "  -0.14%  "

# Row 26 - EthereumClassic
Set-TextValue $ws.Range("E26") "  -2.72%  "

# Row 27 - LidoDAOToken
Set-TextValue $ws.Range("D27") "2.064"
Set-TextValue $ws.Range("E27") "  +0.70%  "

# Row 28 - BitcoinCash
Set-TextValue $ws.Range("D28") "113.12"
Set-TextValue $ws.Range("E28") "  -2.51%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D29") "4.846"
Set-TextValue $ws.Range("E29") "  -3.96%  "

# Row 30 - Stellar
Set-TextValue $ws.Range("D30") "0.08699"
Set-TextValue $ws.Range("E30") "  -1.58%  "

# Row 31 - HuobiToken
Set-TextValue $ws.Range("D31") "3.032"
Set-TextValue $ws.Range("E31") "  -3.57%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "4.507"
Set-TextValue $ws.Range("E32") "  -0.52%  "

# Row 33 - ImmutableX
Set-TextValue $ws.Range("D33") "0.7330"
Set-TextValue $ws.Range("E33") "  -4.96%  "

# Row 34 - RenderToken
Set-TextValue $ws.Range("D34") "2.686"
Set-TextValue $ws.Range("E34") "  -0.20%  "

# Row 35 - ARBITRUM
Set-TextValue $ws.Range("D35") "1.120"
Set-TextValue $ws.Range("E35") "  -4.29%  "

# Row 36 - Frax
Set-TextValue $ws.Range("D36") "1.005"
Set-TextValue $ws.Range("E36") "  +0.54%  "

# Row 37 - TrustWalletToken
Set-TextValue $ws.Range("D37") "1.084"
Set-TextValue $ws.Range("E37") "  -2.40%  "

# Row 38 - VeChain
Set-TextValue $ws.Range("D38") "0.01957"
Set-TextValue $ws.Range("E38") "  +0.61%  "

# Row 39 - Hedera
Set-TextValue $ws.Range("D39") "0.05111"
Set-TextValue $ws.Range("E39") "  -2.37%  "

# Row 40 - MXToken
Set-TextValue $ws.Range("D40") "2.888"
Set-TextValue $ws.Range("E40") "  -2.16%  "

# Row 41 - FraxShare
Set-TextValue $ws.Range("D41") "6.992"
Set-TextValue $ws.Range("E41") "  -1.13%  "

# Row 42 - TheSandbox
Set-TextValue $ws.Range("D42") "0.4995"
Set-TextValue $ws.Range("E42") "  -2.42%  "

# Row 43 - Algorand
Set-TextValue $ws.Range("D43") "0.1555"
Set-TextValue $ws.Range("E43") "  -4.54%  "

# Row 44 - Aptos
Set-TextValue $ws.Range("D44") "8.165"
Set-TextValue $ws.Range("E44") "  -4.27%  "

# Row 45 - PaxDollar
Set-TextValue $ws.Range("D45") "1.008"
Set-TextValue $ws.Range("E45") "  +0.57%  "

# Row 46 - Decentraland
Set-TextValue $ws.Range("D46") "0.4610"
Set-TextValue $ws.Range("E46") "  -4.13%  "

# Row 47 - EnergySwap
Set-TextValue $ws.Range("D47") "10.00"
Set-TextValue $ws.Range("E47") "  -3.42%  "

# Row 48 - Quant
Set-TextValue $ws.Range("D48") "101.16"
Set-TextValue $ws.Range("E48") "  -1.74%  "

# Row 49 - NEARProtocol
Set-TextValue $ws.Range("E49") "  -3.56%  "

# Row 50 - Cronos
Set-TextValue $ws.Range("D50") "0.05996"
Set-TextValue $ws.Range("E50") "  -3.38%  "

# Row 51 - Aave
Set-TextValue $ws.Range("D51") "64.49"
Set-TextValue $ws.Range("E51") "  -1.68%  "
